# Room Information (Responses) - add another form response row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new shared strings in the same order they first appear in the
# target workbook (C5, then D5, then B5) so the sharedStrings table comes
# out in the expected order: "Something, Else", "5:00-5:40 PM", "Some 123.1".
$ws.Range("C5").Value = "Something, Else"
$ws.Range("D5").Value = "5:00-5:40 PM"
$ws.Range("B5").Value = "Some 123.1"

# Timestamp, matching the formatting already used for A2:A4.
$ws.Range("A5").Value = 42070.599336550928
$ws.Range("A5").NumberFormat = "m/d/yyyy\ h:mm:ss"

# Carry the same Tuesday/Wednesday/Thursday selections forward from row 4.
$ws.Range("E5").Value = $ws.Range("E4").Value2
$ws.Range("F5").Value = $ws.Range("F4").Value2
$ws.Range("G5").Value = $ws.Range("G4").Value2

# Match the row height used throughout the rest of the form responses.
$ws.Rows.Item(5).RowHeight = $ws.Rows.Item(4).RowHeight

# Select the newly-entered cell, which also updates the frozen-pane view
# so the sheet no longer scrolls to show column F first.
[void]$ws.Range("B5").Select()
